$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Update Regresi Tanggal" — the date-regression test value in R2 moves
# from 20230505 to 20240907. It's stored as a plain number (cell style
# just renders it as text), so a numeric assignment keeps the same type.
$ws.Range("R2").Value = 20240907

# View state also shifted: the window scrolled one column left (O1 -> N1)
# and the live selection moved down one row, onto the cell that was just
# updated (R2 -> R3).
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("R3").Select()
